$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.912.95"
$ws.Range("E2").Value = "  +2.19%  "
$ws.Range("D3").Value = "3.116.08"
$ws.Range("E3").Value = "  +5.57%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'580.74"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").Value = "'172.97"
$ws.Range("E6").Value = "  +7.50%  "
$ws.Range("D8").Value = "3.110.33"
$ws.Range("E8").Value = "  +5.52%  "
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("E10").Value = "  -3.35%  "
$ws.Range("E11").Value = "  +3.65%  "
$ws.Range("D12").Value = "'0.483"
$ws.Range("E12").Value = "  +4.92%  "
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("E13").Value = "  +2.14%  "
$ws.Range("D14").Value = "'37.26"
$ws.Range("E14").Value = "  +7.65%  "
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "3.631.07"
$ws.Range("E16").Value = "  +5.56%  "
$ws.Range("D17").Value = "66.909.66"
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").Value = "'7.20"
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("D19").Value = "3.114.84"
$ws.Range("E19").Value = "  +5.63%  "
$ws.Range("D20").Value = "'16.21"
$ws.Range("E20").Value = "  +2.47%  "
$ws.Range("D21").Value = "'484.69"
$ws.Range("E21").Value = "  +8.61%  "
$ws.Range("D22").Value = "'0.717"
$ws.Range("E22").Value = "  +3.16%  "
$ws.Range("D23").Value = "'7.54"
$ws.Range("E23").Value = "  +3.28%  "
$ws.Range("D24").Value = "'84.19"
$ws.Range("E24").Value = "  +2.38%  "
$ws.Range("E25").Value = "  +5.51%  "
$ws.Range("D26").Value = "'13.11"
$ws.Range("E26").Value = "  +6.75%  "
$ws.Range("D27").Value = "'10.07"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").Value = "'7.98"
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("E30").Value = "  -3.64%  "
$ws.Range("E31").Value = "  +4.05%  "
$ws.Range("D32").Value = "'29.01"
$ws.Range("E32").Value = "  +6.77%  "
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").Value = "  +2.15%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  +3.84%  "
$ws.Range("D37").Value = "'5.91"
$ws.Range("E37").Value = "  +3.14%  "
$ws.Range("D38").Value = "'48.46"
$ws.Range("E38").Value = "  +7.11%  "
$ws.Range("D39").Value = "'2.12"
$ws.Range("E39").Value = "  +7.44%  "
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("D41").Value = "'0.317"
$ws.Range("E41").Value = "  +5.11%  "
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").Value = "'8.68"
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("E45").Value = "  +3.10%  "
$ws.Range("D46").Value = "2.840.74"
$ws.Range("E46").Value = "  +5.99%  "
$ws.Range("D47").Value = "'381.23"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").Value = "'135.43"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +5.46%  "
$ws.Range("E51").Value = "  +3.27%  "
